$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

# --- settings sheet: upgrade from form_id/form_title style to the newer
#     pyxform/cht-conf settings layout, which drops the form_id column
#     (column B) entirely and shifts version/style/namespaces one column left.

# Capture the comment text from the columns that will shift left before the
# column delete invalidates the column positions.
$versionComment    = $wsSettings.Range("C1").Comment.Text()
$pagesComment      = $wsSettings.Range("D1").Comment.Text()
$namespacesComment = $wsSettings.Range("E1").Comment.Text()

# Remove the form_id column (B); Excel shifts C/D/E left into B/C/D.
$wsSettings.Columns("B").Delete()

# The cell comments themselves are anchored to absolute cell refs and don't
# travel with the shifted cells, so re-home the text manually.
$wsSettings.Range("B1").Comment.Text($versionComment)
$wsSettings.Range("C1").Comment.Text($pagesComment)
$wsSettings.Range("D1").Comment.Text($namespacesComment)
$wsSettings.Range("E1").Comment.Delete()

# --- survey sheet: drop the stray NO_LABEL appearance value left in C3.
$wsSurvey.Range("C3").Clear()

# --- restore a sensible selection/cursor state on both sheets, finishing on
#     the survey sheet so it stays the active tab.
$wsSettings.Range("B1").Select()
$wsSurvey.Range("A11").Select()
